# This workbook gets a weekly update: two new data rows are inserted at the
# top of the "Repollo / Macroferia Regional de Talca" data block (rows 269-270),
# pushing the existing 269-383 block down to 271-385. The two new rows carry a
# fresh date (2022-09-21) and their own price data, and the final dimension
# grows from A1:R383 to A1:R385.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 269; Excel shifts rows 269:383 down to
# 271:385 and copies formatting (e.g. the date style on column D) from the
# row immediately above the insertion point.
$ws.Rows("269:270").Insert()

# --- New row 269 (Primera) ---
$ws.Cells.Item(269, 1).Value  = 5
$ws.Cells.Item(269, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(269, 3).Value  = "Maule"
$ws.Cells.Item(269, 4).Value  = 44825
$ws.Cells.Item(269, 5).Value  = 7
$ws.Cells.Item(269, 6).Value  = 100112006
$ws.Cells.Item(269, 7).Value  = "Repollo"
$ws.Cells.Item(269, 8).Value  = "Crespo record"
$ws.Cells.Item(269, 9).Value  = "Primera"
$ws.Cells.Item(269, 10).Value = 2000
$ws.Cells.Item(269, 11).Value = 1300
$ws.Cells.Item(269, 12).Value = 1300
$ws.Cells.Item(269, 13).Value = 1300
$ws.Cells.Item(269, 14).Value = "`$/unidad"
$ws.Cells.Item(269, 15).Value = "Región del Maule"
$ws.Cells.Item(269, 16).Value = 1300
$ws.Cells.Item(269, 17).Value = 1
$ws.Cells.Item(269, 18).Value = "Hortaliza"

# --- New row 270 (Segunda) ---
$ws.Cells.Item(270, 1).Value  = 5
$ws.Cells.Item(270, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(270, 3).Value  = "Maule"
$ws.Cells.Item(270, 4).Value  = 44825
$ws.Cells.Item(270, 5).Value  = 7
$ws.Cells.Item(270, 6).Value  = 100112006
$ws.Cells.Item(270, 7).Value  = "Repollo"
$ws.Cells.Item(270, 8).Value  = "Crespo record"
$ws.Cells.Item(270, 9).Value  = "Segunda"
$ws.Cells.Item(270, 10).Value = 2000
$ws.Cells.Item(270, 11).Value = 1000
$ws.Cells.Item(270, 12).Value = 1000
$ws.Cells.Item(270, 13).Value = 1000
$ws.Cells.Item(270, 14).Value = "`$/unidad"
$ws.Cells.Item(270, 15).Value = "Región del Maule"
$ws.Cells.Item(270, 16).Value = 1000
$ws.Cells.Item(270, 17).Value = 1
$ws.Cells.Item(270, 18).Value = "Hortaliza"
